$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Update Sheet2 header row: measurement-level / data-type labels
$ws2.Range("B1").Value = "Double"
$ws2.Range("C1").Value = "string"
$ws2.Range("D1").Value = "integer"

# Update Sheet2 data rows - B and D columns become numeric, C stays text
$ws2.Range("B2").Value = 0.76
$ws2.Range("D2").Value = 1700

$ws2.Range("B3").Value = 1.2
$ws2.Range("D3").Value = 12

# Update selections on each sheet
$null = $ws1.Range("E16").Select()
$null = $ws2.Range("C22").Select()

# Make Sheet2 the active tab
$null = $ws2.Activate()
